# Apply updated cryptos list values (prices, volumes, and a couple of row reorderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to Text format first so numeric-looking strings
# (e.g. "26.80", "103.00", "30.212.46") are preserved verbatim instead of
# being coerced into floating point numbers by Excel.
$priceCells = @("D2","D3","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D16","D17","D18","D19","D20","D21","D23","D24","D25","D26","D27","D28","D29","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.212.46"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").Value = "1.913.38"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "0.8194"
$ws.Range("E5").Value = "  +3.68%  "

$ws.Range("D6").Value = "243.67"
$ws.Range("E6").Value = "  +0.32%  "

$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "0.3254"
$ws.Range("E8").Value = "  +2.80%  "

$ws.Range("D9").Value = "26.80"

$ws.Range("D10").Value = "0.07067"
$ws.Range("E10").Value = "  +2.06%  "

$ws.Range("D11").Value = "0.08099"
$ws.Range("E11").Value = "  +1.46%  "

$ws.Range("D12").Value = "0.7764"
$ws.Range("E12").Value = "  +3.75%  "

$ws.Range("D13").Value = "1.918.21"
$ws.Range("E13").Value = "  +0.57%  "

$ws.Range("D14").Value = "5.301"
$ws.Range("E14").Value = "  +1.26%  "

$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").Value = "30.216.31"
$ws.Range("E16").Value = "  +0.50%  "

$ws.Range("D17").Value = "14.26"
$ws.Range("E17").Value = "  +1.58%  "

$ws.Range("D18").Value = "5.929"
$ws.Range("E18").Value = "  -0.26%  "

$ws.Range("D19").Value = "247.35"
$ws.Range("E19").Value = "  +0.21%  "

$ws.Range("D20").Value = "0.000007795"

$ws.Range("D21").Value = "2.166.36"
$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D23").Value = "1.001"

$ws.Range("D24").Value = "7.155"
$ws.Range("E24").Value = "  +3.47%  "

$ws.Range("D25").Value = "0.1677"
$ws.Range("E25").Value = "  +20.61%  "

$ws.Range("D26").Value = "9.342"
$ws.Range("E26").Value = "  +0.29%  "

$ws.Range("D27").Value = "167.57"
$ws.Range("E27").Value = "  -1.38%  "

$ws.Range("D28").Value = "19.00"
$ws.Range("E28").Value = "  +0.36%  "

$ws.Range("D29").Value = "2.112"
$ws.Range("E29").Value = "  +3.29%  "

$ws.Range("E30").Value = "  -0.33%  "

$ws.Range("D31").Value = "1.530"
$ws.Range("E31").Value = "  +0.46%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "4.314"
$ws.Range("E32").Value = "  -0.70%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.05833"
$ws.Range("E33").Value = "  +4.42%  "

$ws.Range("D34").Value = "4.107"
$ws.Range("E34").Value = "  -0.18%  "

$ws.Range("D35").Value = "1.278"
$ws.Range("E35").Value = "  +1.07%  "

$ws.Range("D36").Value = "0.7381"
$ws.Range("E36").Value = "  +0.13%  "

$ws.Range("D37").Value = "2.705"
$ws.Range("E37").Value = "  -0.90%  "

$ws.Range("D38").Value = "0.01928"
$ws.Range("E38").Value = "  -0.89%  "

$ws.Range("D39").Value = "2.802"
$ws.Range("E39").Value = "  +0.32%  "

$ws.Range("D40").Value = "0.4471"
$ws.Range("E40").Value = "  +0.58%  "

$ws.Range("D41").Value = "73.44"
$ws.Range("E41").Value = "  +1.09%  "

$ws.Range("D42").Value = "5.976"
$ws.Range("E42").Value = "  -3.45%  "

$ws.Range("D43").Value = "0.8544"
$ws.Range("E43").Value = "  +2.29%  "

$ws.Range("D44").Value = "1.915"
$ws.Range("E44").Value = "  +0.77%  "

$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").Value = "103.00"
$ws.Range("E46").Value = "  +2.40%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.010.09"
$ws.Range("E47").Value = "  +2.30%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "7.595"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "9.866"
$ws.Range("E49").Value = "  +0.58%  "

$ws.Range("D50").Value = "2.068.08"
$ws.Range("E50").Value = "  +0.35%  "

$ws.Range("D51").Value = "1.568"
$ws.Range("E51").Value = "  +4.46%  "
